# Append: 2025-09-17 06:33 JST
# Update the "取得日時" (acquisition timestamp) column for all existing
# data rows on the "ランサーズ" sheet from 06:26:41 to 06:33:15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-09-17 06:26:41"
$newTimestamp = "2025-09-17 06:33:15"

for ($row = 2; $row -le 16; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
